# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# This updates the "K" column (column G) values for the game log rows on
# Sheet1. The values below are the newly (re)computed K (strikeouts)
# figures that replace the previous placeholder/Strike# derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = [ordered]@{
    2  = 2
    3  = 2
    4  = 0
    5  = 2
    6  = 2
    7  = 1
    8  = 2
    9  = 1
    10 = 0
    11 = 4
    13 = 1
    14 = 1
    16 = 2
    17 = 0
    18 = 2
    19 = 3
    21 = 3
    22 = 0
    23 = 2
    24 = 0
    25 = 1
    26 = 0
    27 = 1
    28 = 2
    29 = 0
    30 = 0
    31 = 1
    32 = 0
    33 = 0
    34 = 0
    35 = 1
    36 = 0
    37 = 3
    38 = 2
    39 = 0
    41 = 1
    42 = 1
    43 = 2
    45 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
